$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.235.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = "'1.830.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = "'236.80"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'0.6070"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.50%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = "'0.07120"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.77%  '
$ws.Range("D9").Value = "'0.2819"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").Value = "'23.94"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").Value = "'0.07670"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = "'1.818.06"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").Value = "'4.820"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.23%  '
$ws.Range("D14").Value = "'0.6368"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.92%  '
$ws.Range("D15").Value = "'0.00001001"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").Value = "'2.082.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").Value = "'79.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("D18").Value = "'5.915"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -5.09%  '
$ws.Range("D19").Value = "'29.226.24"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = "'228.93"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = "'11.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.97%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = "'7.015"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.25%  '
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = "'154.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("D26").Value = "'8.106"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.1288"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.89%  '
$ws.Range("D28").Value = "'16.66"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.71%  '
$ws.Range("D29").Value = "'1.498"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.63%  '
$ws.Range("D30").Value = "'0.06476"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("E31").Value = '  -2.21%  '
$ws.Range("D32").Value = "'3.842"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.45%  '
$ws.Range("D33").Value = "'3.845"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.41%  '
$ws.Range("D34").Value = "'1.131"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = "'1.740"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.20%  '
$ws.Range("D36").Value = "'0.6548"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.39%  '
$ws.Range("D37").Value = "'2.540"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("D38").Value = "'2.766"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("D39").Value = "'1.223.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").Value = "'0.01757"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.92%  '
$ws.Range("D41").Value = "'6.590"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.56%  '
$ws.Range("D42").Value = "'0.9260"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = "'101.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = "'1.983.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").Value = "'1.609"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.20%  '
$ws.Range("D49").Value = "'8.541"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.24%  '
$ws.Range("D50").Value = "'6.506"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -7.70%  '
$ws.Range("D51").Value = "'0.05543"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.35%  '
